$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(15, 1).Value = 'мелочь'
$ws.Cells.Item(16, 1).Value = 'особливый товар'
$ws.Cells.Item(17, 1).Value = 'серебреный товар'
$ws.Cells.Item(18, 1).Value = 'деревенский товар'
$ws.Cells.Item(19, 1).Value = 'небогатый товар'
$ws.Cells.Item(20, 1).Value = 'крамными товар'
$ws.Cells.Item(21, 1).Value = 'железный товар'
$ws.Cells.Item(22, 1).Value = 'мясо'
$ws.Cells.Item(24, 1).Value = 'набойчатый товар'
$ws.Cells.Item(25, 1).Value = 'нужный товар'
$ws.Cells.Item(26, 1).Value = 'щепетильный товар'
$ws.Cells.Item(27, 1).Value = 'пушной товар'
$ws.Cells.Item(28, 1).Value = 'суровский товар'
$ws.Cells.Item(30, 1).Value = 'медный товар'
$ws.Cells.Item(31, 1).Value = 'внутренний товар'
$ws.Cells.Item(32, 1).Value = 'питейный припасы'
$ws.Cells.Item(35, 1).Value = 'галантерейный товар'
$ws.Cells.Item(38, 1).Value = 'купецкий товар'
$ws.Cells.Item(39, 1).Value = 'домовый товар'
$ws.Cells.Item(40, 1).Value = 'рукодельный товар'
$ws.Cells.Item(41, 1).Value = 'надлежащий товар'
$ws.Cells.Item(42, 1).Value = 'харчевой припасы'
$ws.Cells.Item(43, 1).Value = 'меховой товар'
